$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 10396.091
$ws.Range("I6").Value = 1206.3334
$ws.Range("J6").Value = 51750
$ws.Range("K6").Value = 3619.0002
$ws.Range("L6").Value = 155250
$ws.Range("M6").Value = -3507.0002
$ws.Range("N6").Value = -155474
$ws.Range("H17").Value = 11114521
$ws.Range("J17").Value = 11631436
$ws.Range("L17").Value = 34894308
$ws.Range("N17").Value = -34894644
$ws.Range("H39").Value = 331.25
$ws.Range("I39").Value = 62.5
$ws.Range("J39").Value = 600
$ws.Range("K39").Value = 187.5
$ws.Range("L39").Value = 1800
$ws.Range("M39").Value = 108.5
$ws.Range("N39").Value = -2392
$ws.Range("H62").Value = 5625.125
$ws.Range("I62").Value = 8998.75
$ws.Range("J62").Value = 2251.5
$ws.Range("K62").Value = 8998.75
$ws.Range("L62").Value = 2251.5
$ws.Range("M62").Value = -8374.75
$ws.Range("N62").Value = -3499.5
$ws.Range("H65").Value = 5625.125
$ws.Range("I65").Value = 8998.75
$ws.Range("J65").Value = 2251.5
$ws.Range("K65").Value = 44993.75
$ws.Range("L65").Value = 11257.5
$ws.Range("M65").Value = -41873.75
$ws.Range("N65").Value = -17497.5
$ws.Range("H94").Value = 4579.9165
$ws.Range("I94").Value = 4579.9165
$ws.Range("K94").Value = 4579.9165
$ws.Range("M94").Value = -4128.9165
$ws.Range("H97").Value = 12188
$ws.Range("J97").Value = 12188
$ws.Range("L97").Value = 36564
$ws.Range("N97").Value = -37556
$ws.Range("H99").Value = 1683.8572
$ws.Range("J99").Value = 150
$ws.Range("L99").Value = 450
$ws.Range("N99").Value = -3446
$ws.Range("H100").Value = 1497.8235
$ws.Range("I100").Value = 1618.75
$ws.Range("J100").Value = 1207.6
$ws.Range("K100").Value = 1618.75
$ws.Range("L100").Value = 1207.6
$ws.Range("M100").Value = -1077.75
$ws.Range("N100").Value = -2289.6
$ws.Range("H101").Value = 684.6667
$ws.Range("I101").Value = 465
$ws.Range("J101").Value = 794.5
$ws.Range("K101").Value = 1395
$ws.Range("L101").Value = 2383.5
$ws.Range("M101").Value = 227
$ws.Range("N101").Value = -5627.5
$ws.Range("H129").Value = 2868.15
$ws.Range("J129").Value = 2592.5557
$ws.Range("L129").Value = 7777.6671
$ws.Range("N129").Value = -17777.6671
$ws.Range("H131").Value = 2990.35
$ws.Range("I131").Value = 2930
$ws.Range("J131").Value = 3005.4375
$ws.Range("K131").Value = 8790
$ws.Range("L131").Value = 9016.3125
$ws.Range("M131").Value = -3750
$ws.Range("N131").Value = -19096.3125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1144.591
$ws.Range("I97").Value = 984.1667
$ws.Range("J97").Value = 1337.1
$ws.Range("K97").Value = 984.1667
$ws.Range("L97").Value = 1337.1
$ws.Range("M97").Value = -488.1667
$ws.Range("N97").Value = -2329.1
$ws.Range("H122").Value = 2045.2122
$ws.Range("I122").Value = 2147.5833
$ws.Range("J122").Value = 1772.2222
$ws.Range("K122").Value = 6442.749899999999
$ws.Range("L122").Value = 5316.6666
$ws.Range("M122").Value = -3992.749899999999
$ws.Range("N122").Value = -10216.6666
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 570.6070999999999
$ws.Range("I94").Value = 590.86957
$ws.Range("J94").Value = 477.4
$ws.Range("K94").Value = 590.86957
$ws.Range("L94").Value = 477.4
$ws.Range("M94").Value = -139.86957
$ws.Range("N94").Value = -1379.4
$ws.Range("H99").Value = 1947.2
$ws.Range("I99").Value = 1870.8182
$ws.Range("J99").Value = 2507.3333
$ws.Range("K99").Value = 1870.8182
$ws.Range("L99").Value = 2507.3333
$ws.Range("M99").Value = -372.8181999999999
$ws.Range("N99").Value = -5503.3333
$ws.Range("H106").Value = 44100.668
$ws.Range("J106").Value = 44100.668
$ws.Range("L106").Value = 44100.668
$ws.Range("N106").Value = -46624.668
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 47322
$ws.Range("J43").Value = 47322
$ws.Range("L43").Value = 47322
$ws.Range("N43").Value = -47690
$ws.Range("H101").Value = 47322
$ws.Range("J101").Value = 47322
$ws.Range("L101").Value = 47322
$ws.Range("N101").Value = -53812
$ws.Range("H105").Value = 3390.625
$ws.Range("I105").Value = 3473.3333
$ws.Range("K105").Value = 3473.3333
$ws.Range("M105").Value = -1726.3333
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3189.7896
$ws.Range("I64").Value = 1435.1666
$ws.Range("J64").Value = 3999.6155
$ws.Range("K64").Value = 4305.4998
$ws.Range("L64").Value = 11998.8465
$ws.Range("M64").Value = -4035.4998
$ws.Range("N64").Value = -12538.8465
$ws.Range("H67").Value = 3189.7896
$ws.Range("I67").Value = 1435.1666
$ws.Range("J67").Value = 3999.6155
$ws.Range("K67").Value = 4305.4998
$ws.Range("L67").Value = 11998.8465
$ws.Range("M67").Value = -3369.4998
$ws.Range("N67").Value = -13870.8465
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 996.9375
$ws.Range("I122").Value = 1050.1538
$ws.Range("J122").Value = 766.3333
$ws.Range("K122").Value = 3150.4614
$ws.Range("L122").Value = 2298.9999
$ws.Range("M122").Value = -700.4614000000001
$ws.Range("N122").Value = -7198.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4893.1763
$ws.Range("I46").Value = 677.5
$ws.Range("J46").Value = 6190.3076
$ws.Range("K46").Value = 677.5
$ws.Range("L46").Value = 6190.3076
$ws.Range("M46").Value = -489.5
$ws.Range("N46").Value = -6566.3076
$ws.Range("H61").Value = 3339.4614
$ws.Range("I61").Value = 3663.5
$ws.Range("J61").Value = 2821
$ws.Range("K61").Value = 3663.5
$ws.Range("L61").Value = 2821
$ws.Range("M61").Value = -3461.5
$ws.Range("N61").Value = -3225
$ws.Range("H113").Value = 3339.4614
$ws.Range("I113").Value = 3663.5
$ws.Range("J113").Value = 2821
$ws.Range("K113").Value = 3663.5
$ws.Range("L113").Value = 2821
$ws.Range("M113").Value = -1493.5
$ws.Range("N113").Value = -7161
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4340
$ws.Range("I62").Value = 2700
$ws.Range("J62").Value = 4750
$ws.Range("K62").Value = 2700
$ws.Range("L62").Value = 4750
$ws.Range("M62").Value = -2076
$ws.Range("N62").Value = -5998
$ws.Range("H65").Value = 4340
$ws.Range("I65").Value = 2700
$ws.Range("J65").Value = 4750
$ws.Range("K65").Value = 13500
$ws.Range("L65").Value = 23750
$ws.Range("M65").Value = -10380
$ws.Range("N65").Value = -29990
$ws.Range("H96").Value = 1475.3334
$ws.Range("J96").Value = 1796.8334
$ws.Range("L96").Value = 1796.8334
$ws.Range("N96").Value = -4542.8334
$ws.Range("H100").Value = 550
$ws.Range("I100").Value = 550
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1100
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -559
$ws.Range("N100").ClearContents()
